$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C16").Value = "1050957805"
$ws.Range("D16").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E16").Value = "2104"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 877803

$ws.Range("C17").Value = "1050957805"
$ws.Range("D17").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E17").Value = "2105"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 877803

$ws.Range("C18").Value = "1050957805"
$ws.Range("D18").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E18").Value = "2106"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 877803

$ws.Range("C19").Value = "1050957805"
$ws.Range("D19").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E19").Value = "2107"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 877803

$ws.Range("C20").Value = "1050957805"
$ws.Range("D20").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E20").Value = "2108"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 877803

$ws.Range("C21").Value = "1050957805"
$ws.Range("D21").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E21").Value = "2109"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 877803

$ws.Range("C22").Value = "1050957805"
$ws.Range("D22").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E22").Value = "2110"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 877803

$ws.Range("C23").Value = "1050957805"
$ws.Range("D23").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E23").Value = "2111"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 877803

$ws.Range("C24").Value = "1196966238"
$ws.Range("D24").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E24").Value = "2111"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 877803

$ws.Range("C25").Value = "1050957805"
$ws.Range("D25").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E25").Value = "2112"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 877803

$ws.Range("C26").Value = "1196966238"
$ws.Range("D26").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E26").Value = "2112"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 877803

$ws.Range("C27").Value = "1050957805"
$ws.Range("D27").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E27").Value = "2201"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 877803

$ws.Range("C28").Value = "1196966238"
$ws.Range("D28").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E28").Value = "2201"
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 877803

$ws.Range("C29").Value = "1050957805"
$ws.Range("D29").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E29").Value = "2202"
$ws.Range("F29").Value = 35112
$ws.Range("G29").Value = 877803

$ws.Range("C30").Value = "1196966238"
$ws.Range("D30").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E30").Value = "2202"
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 877803

$ws.Range("C31").Value = "1050957805"
$ws.Range("D31").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E31").Value = "2203"
$ws.Range("F31").Value = 35112
$ws.Range("G31").Value = 877803

$ws.Range("C32").Value = "1196966238"
$ws.Range("D32").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E32").Value = "2203"
$ws.Range("F32").Value = 36341
$ws.Range("G32").Value = 877803

$ws.Range("C33").Value = "1050957805"
$ws.Range("D33").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E33").Value = "2204"
$ws.Range("F33").Value = 35112
$ws.Range("G33").Value = 877803

$ws.Range("C34").Value = "1196966238"
$ws.Range("D34").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E34").Value = "2204"
$ws.Range("F34").Value = 36341
$ws.Range("G34").Value = 877803

$ws.Range("C35").Value = "1050957805"
$ws.Range("D35").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E35").Value = "2205"
$ws.Range("F35").Value = 35112
$ws.Range("G35").Value = 877803

$ws.Range("C36").Value = "1196966238"
$ws.Range("D36").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E36").Value = "2205"
$ws.Range("F36").Value = 36341
$ws.Range("G36").Value = 877803

$ws.Range("C37").Value = "1050957805"
$ws.Range("D37").Value = "ANGEL MARTINEZ VASQUEZ MARTINEZ"
$ws.Range("E37").Value = "2206"
$ws.Range("F37").Value = 25749
$ws.Range("G37").Value = 877803

$ws.Range("C38").Value = "1196966238"
$ws.Range("D38").Value = "CARLOS ALFREDO ROSARIO PENATE"
$ws.Range("E38").Value = "2206"
$ws.Range("F38").Value = 25749
$ws.Range("G38").Value = 877803
